# Applies the commit:
#  1) Updated the selector (Login_Info!A2/B2 now hold the real
#     login e-mail / password, A2 is a mailto hyperlink rendered with the
#     Webdings icon font; Friend_Request_Management becomes the active tab
#     with a fresh selection) and
#  2) Added the Delay in the login button click activity (Friend_Request_Management!A2
#     flips from 9 -> 1 and a new Delay-in-seconds value of 50 is written to B2).

$wb  = $excel.ActiveWorkbook
$wsLogin = $wb.Worksheets.Item("Login_Info")
$wsFrm   = $wb.Worksheets.Item("Friend_Request_Management")

# ---- Login_Info sheet ---------------------------------------------------
# Row 2 grows a little taller to match the new content.
$wsLogin.Rows.Item(2).RowHeight = 15.75

# A2: the e-mail address, written out as a clickable mailto hyperlink that
# renders using the Webdings icon font (mirrors the Hyperlink cell style).
$wsLogin.Range("A2").Value = "rhassnain@gmail.com"
$wsLogin.Hyperlinks.Add($wsLogin.Range("A2"), "mailto:rhassnain@gmail.com") | Out-Null
$wsLogin.Range("A2").Font.Name = "Webdings"
$wsLogin.Range("A2").Font.Underline = $true

# B2: the password, rendered in the Webdings font (not a hyperlink).
$wsLogin.Range("B2").Value = "Raza#123"
$wsLogin.Range("B2").Font.Name = "Webdings"

$wsLogin.PageSetup.PaperSize = 9
$wsLogin.PageSetup.Orientation = 1

# ---- Friend_Request_Management sheet -----------------------------------
# A2: flip the flag from 9 to 1.
$wsFrm.Range("A2").Value = 1
# B2: the new Delay (seconds) used by the login button click activity.
$wsFrm.Range("B2").Value = 50

# ---- Selections / active tab --------------------------------------------
$wsLogin.Range("B21").Select() | Out-Null
$wsFrm.Activate() | Out-Null
$wsFrm.Range("B8").Select() | Out-Null
